$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (rows 6 and 7) so only 4 data rows remain (rows 2-5)
$ws.Rows(6).Delete()
$ws.Rows(6).Delete()

# Row 2: FAPs -> FAPs
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Rspo3"
$ws.Range("C2").Value2 = "Lgr5"
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 2.507621333333333
$ws.Range("H2").Value2 = 7.522864
$ws.Range("I2").Value2 = 0.9300694554254023
$ws.Range("J2").Value2 = 0.9300694554254023
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.8326193333333333
$ws.Range("N2").Value2 = 2.497858
$ws.Range("O2").Value2 = 0.9388124812781204
$ws.Range("P2").Value2 = 0.9388124812781203
$ws.Range("Q2").Value2 = 2.087894002812444
$ws.Range("R2").Value2 = 18.791046025312
$ws.Range("S2").Value2 = 0.8731608132089121
$ws.Range("T2").Value2 = 0.873160813208912

# Row 3: FAPs -> MuSCs
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Rspo3"
$ws.Range("C3").Value2 = "Lgr5"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 2.507621333333333
$ws.Range("H3").Value2 = 7.522864
$ws.Range("I3").Value2 = 0.9300694554254023
$ws.Range("J3").Value2 = 0.9300694554254023
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.05426633333333333
$ws.Range("N3").Value2 = 0.162799
$ws.Range("O3").Value2 = 0.0611875187218796
$ws.Range("P3").Value2 = 0.06118751872187959
$ws.Range("Q3").Value2 = 0.1360794151484444
$ws.Range("R3").Value2 = 1.224714736336
$ws.Range("S3").Value2 = 0.05690864221649017
$ws.Range("T3").Value2 = 0.05690864221649016

# Row 4: MuSCs -> FAPs
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("B4").Value2 = "Rspo3"
$ws.Range("C4").Value2 = "Lgr5"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.1885443333333333
$ws.Range("H4").Value2 = 0.5656329999999999
$ws.Range("I4").Value2 = 0.06993054457459773
$ws.Range("J4").Value2 = 0.06993054457459771
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.8326193333333333
$ws.Range("N4").Value2 = 2.497858
$ws.Range("O4").Value2 = 0.9388124812781204
$ws.Range("P4").Value2 = 0.9388124812781203
$ws.Range("Q4").Value2 = 0.1569856571237777
$ws.Range("R4").Value2 = 1.412870914114
$ws.Range("S4").Value2 = 0.0656516680692083
$ws.Range("T4").Value2 = 0.06565166806920827

# Row 5: MuSCs -> MuSCs
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Rspo3"
$ws.Range("C5").Value2 = "Lgr5"
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.1885443333333333
$ws.Range("H5").Value2 = 0.5656329999999999
$ws.Range("I5").Value2 = 0.06993054457459773
$ws.Range("J5").Value2 = 0.06993054457459771
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.05426633333333333
$ws.Range("N5").Value2 = 0.162799
$ws.Range("O5").Value2 = 0.0611875187218796
$ws.Range("P5").Value2 = 0.06118751872187959
$ws.Range("Q5").Value2 = 0.01023160964077778
$ws.Range("R5").Value2 = 0.09208448676699998
$ws.Range("S5").Value2 = 0.004278876505389434
$ws.Range("T5").Value2 = 0.004278876505389433
